# "Monday and Tuesday results" - two pairs of whole sample records were
# entered on the wrong row. Fix by swapping every column (including A,
# sampleid) between row 7 <-> row 8, and between row 13 <-> row 14.
#
# Use Copy/Paste (not .Value2 array assignment) so each cell keeps its
# original type/format (e.g. numeric-looking sample ids like "424" stay
# text, dates stay dates) exactly as it was on the row it moves to - a
# staging row far below the data is used to do a true 2-way swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 7 and 8 in full (sampleid 280 and 424 trade places) ---
$ws.Range("A7:CX7").Copy($ws.Range("A1000:CX1000"))
$ws.Range("A8:CX8").Copy($ws.Range("A7:CX7"))
$ws.Range("A1000:CX1000").Copy($ws.Range("A8:CX8"))
$ws.Range("A1000:CX1000").Clear()

# --- Swap rows 13 and 14 in full (sampleid 299 and 253 trade places) ---
$ws.Range("A13:CX13").Copy($ws.Range("A1000:CX1000"))
$ws.Range("A14:CX14").Copy($ws.Range("A13:CX13"))
$ws.Range("A1000:CX1000").Copy($ws.Range("A14:CX14"))
$ws.Range("A1000:CX1000").Clear()
